# Update "想去人数" (F column) counts on multiple sheets to reflect the
# newly generated data snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 644
$ws1.Range("F8").Value = 355
$ws1.Range("F9").Value = 1772
$ws1.Range("F11").Value = 1444
$ws1.Range("F13").Value = 354
$ws1.Range("F14").Value = 697
$ws1.Range("F15").Value = 12945
$ws1.Range("F16").Value = 12892
$ws1.Range("F21").Value = 58
$ws1.Range("F22").Value = 600
$ws1.Range("F23").Value = 2026
$ws1.Range("F24").Value = 39
$ws1.Range("F25").Value = 18
$ws1.Range("F26").Value = 12
$ws1.Range("F28").Value = 114
$ws1.Range("F30").Value = 698

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 20

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F11").Value = 644
$ws4.Range("F13").Value = 355
$ws4.Range("F14").Value = 1772
$ws4.Range("F16").Value = 1444
$ws4.Range("F18").Value = 354
$ws4.Range("F20").Value = 697
$ws4.Range("F21").Value = 12945
$ws4.Range("F22").Value = 12892
$ws4.Range("F27").Value = 58
$ws4.Range("F28").Value = 600
$ws4.Range("F29").Value = 20
$ws4.Range("F31").Value = 2026
$ws4.Range("F32").Value = 39
$ws4.Range("F33").Value = 18
$ws4.Range("F34").Value = 12
$ws4.Range("F38").Value = 114
$ws4.Range("F40").Value = 698
